# Insert 3 new rows before row 1174, shifting existing rows 1174:1234 down to 1177:1237,
# then populate the 3 new rows with their data (weekly update: new records for date 44753,
# and earlier "Larga vida" duplicate records renumbered down the sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at 1174, pushing the rest of the table down by 3.
$ws.Rows("1174:1176").Insert()

# Common (unchanged across these rows) column values for this data block.
$A = 8
$B = "Terminal La Palmera de La Serena"
$C = "Coquimbo"
$E = 4
$F = 100112020
$G = "Tomate"
$R = "Hortaliza"

# --- New row 1174 ---
$ws.Cells.Item(1174, 1).Value  = $A
$ws.Cells.Item(1174, 2).Value  = $B
$ws.Cells.Item(1174, 3).Value  = $C
$ws.Cells.Item(1174, 4).Value  = 44753
$ws.Cells.Item(1174, 5).Value  = $E
$ws.Cells.Item(1174, 6).Value  = $F
$ws.Cells.Item(1174, 7).Value  = $G
$ws.Cells.Item(1174, 8).Value  = "Larga vida"
$ws.Cells.Item(1174, 9).Value  = "Primera"
$ws.Cells.Item(1174, 10).Value = 400
$ws.Cells.Item(1174, 11).Value = 8000
$ws.Cells.Item(1174, 12).Value = 8500
$ws.Cells.Item(1174, 13).Value = 8250
$ws.Cells.Item(1174, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(1174, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1174, 16).Value = 458
$ws.Cells.Item(1174, 17).Value = 18
$ws.Cells.Item(1174, 18).Value = $R

# --- New row 1175 ---
$ws.Cells.Item(1175, 1).Value  = $A
$ws.Cells.Item(1175, 2).Value  = $B
$ws.Cells.Item(1175, 3).Value  = $C
$ws.Cells.Item(1175, 4).Value  = 44753
$ws.Cells.Item(1175, 5).Value  = $E
$ws.Cells.Item(1175, 6).Value  = $F
$ws.Cells.Item(1175, 7).Value  = $G
$ws.Cells.Item(1175, 8).Value  = "Larga vida"
$ws.Cells.Item(1175, 9).Value  = "Primera"
$ws.Cells.Item(1175, 10).Value = 600
$ws.Cells.Item(1175, 11).Value = 4800
$ws.Cells.Item(1175, 12).Value = 5000
$ws.Cells.Item(1175, 13).Value = 4900
$ws.Cells.Item(1175, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(1175, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1175, 16).Value = 490
$ws.Cells.Item(1175, 17).Value = 10
$ws.Cells.Item(1175, 18).Value = $R

# --- New row 1176 ---
$ws.Cells.Item(1176, 1).Value  = $A
$ws.Cells.Item(1176, 2).Value  = $B
$ws.Cells.Item(1176, 3).Value  = $C
$ws.Cells.Item(1176, 4).Value  = 44753
$ws.Cells.Item(1176, 5).Value  = $E
$ws.Cells.Item(1176, 6).Value  = $F
$ws.Cells.Item(1176, 7).Value  = $G
$ws.Cells.Item(1176, 8).Value  = "Larga vida"
$ws.Cells.Item(1176, 9).Value  = "Segunda"
$ws.Cells.Item(1176, 10).Value = 400
$ws.Cells.Item(1176, 11).Value = 3800
$ws.Cells.Item(1176, 12).Value = 4000
$ws.Cells.Item(1176, 13).Value = 3900
$ws.Cells.Item(1176, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(1176, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(1176, 16).Value = 390
$ws.Cells.Item(1176, 17).Value = 10
$ws.Cells.Item(1176, 18).Value = $R
